$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D8").Value = 52838900
$ws.Range("E8").Value = 50448300
$ws.Range("F8").Value = 46224100
$ws.Range("G8").Value = 43872100
$ws.Range("H8").Value = 40657100
$ws.Range("I8").Value = 40086600
$ws.Range("J8").Value = 35685100
$ws.Range("D9").Value = 25670500
$ws.Range("E9").Value = 25087300
$ws.Range("F9").Value = 21736100
$ws.Range("G9").Value = 19969200
$ws.Range("H9").Value = 18530800
$ws.Range("I9").Value = 17642900
$ws.Range("J9").Value = 14977800
$ws.Range("D10").Value = 27168400
$ws.Range("E10").Value = 25361000
$ws.Range("F10").Value = 24488100
$ws.Range("G10").Value = 23902900
$ws.Range("H10").Value = 22126300
$ws.Range("I10").Value = 22443700
$ws.Range("J10").Value = 20707300
$ws.Range("D15").Value = 8284200
$ws.Range("E15").Value = 7681800
$ws.Range("F15").Value = 6502000
$ws.Range("G15").Value = 5947500
$ws.Range("H15").Value = 5251400
$ws.Range("I15").Value = 5357400
$ws.Range("J15").Value = 4861500
$ws.Range("D17").Value = 47659500
$ws.Range("E17").Value = 44779300
$ws.Range("F17").Value = 38910200
$ws.Range("G17").Value = 35775100
$ws.Range("H17").Value = 32678900
$ws.Range("I17").Value = 31751900
$ws.Range("J17").Value = 27533700
$ws.Range("D18").Value = 5179400
$ws.Range("E18").Value = 5669000
$ws.Range("F18").Value = 7313900
$ws.Range("G18").Value = 8097000
$ws.Range("H18").Value = 7978200
$ws.Range("I18").Value = 8334700
$ws.Range("J18").Value = 8151400
$ws.Range("D20").Value = -659200
$ws.Range("E20").Value = -2703400
$ws.Range("F20").Value = -1183100
$ws.Range("G20").Value = -1956400
$ws.Range("H20").Value = -1289900
$ws.Range("I20").Value = 72300
$ws.Range("J20").Value = -458700
$ws.Range("D21").Value = 12796400
$ws.Range("E21").Value = 10639900
$ws.Range("F21").Value = 12626500
$ws.Range("G21").Value = 12082300
$ws.Range("H21").Value = 11934600
$ws.Range("I21").Value = 13759200
$ws.Range("J21").Value = "NA"
$ws.Range("D22").Value = 1567200
$ws.Range("E22").Value = 1751300
$ws.Range("F22").Value = 3227200
$ws.Range("G22").Value = 1630300
$ws.Range("H22").Value = 1238700
$ws.Range("I22").Value = 1288600
$ws.Range("J22").Value = 1075300
$ws.Range("D23").Value = 2953000
$ws.Range("E23").Value = 1214300
$ws.Range("F23").Value = 2903600
$ws.Range("G23").Value = 4510300
$ws.Range("H23").Value = 5449600
$ws.Range("I23").Value = 7118400
$ws.Range("J23").Value = 6617300
$ws.Range("D24").Value = 1290000
$ws.Range("E24").Value = 589500
$ws.Range("F24").Value = 992000
$ws.Range("G24").Value = 2053700
$ws.Range("H24").Value = 1571900
$ws.Range("I24").Value = 2378300
$ws.Range("J24").Value = 2055700
$ws.Range("D26").Value = 1663100
$ws.Range("E26").Value = 624700
$ws.Range("F26").Value = 1911600
$ws.Range("G26").Value = 2456600
$ws.Range("H26").Value = 3877700
$ws.Range("I26").Value = 4740100
$ws.Range("J26").Value = 4561700
$ws.Range("D27").Value = 1516700
$ws.Range("E27").Value = 447300
$ws.Range("F27").Value = 1813000
$ws.Range("G27").Value = 2386700
$ws.Range("H27").Value = 3859600
$ws.Range("I27").Value = 4705900
$ws.Range("J27").Value = 4295100
$ws.Range("D32").Value = 659200
$ws.Range("E32").Value = 2703400
$ws.Range("F32").Value = 1183100
$ws.Range("G32").Value = 1956400
$ws.Range("H32").Value = 1289900
$ws.Range("I32").Value = -72300
$ws.Range("J32").Value = 458700
$ws.Range("D33").Value = 1516700
$ws.Range("E33").Value = 447300
$ws.Range("F33").Value = 1813000
$ws.Range("G33").Value = 2386700
$ws.Range("H33").Value = 3859600
$ws.Range("I33").Value = 4705900
$ws.Range("J33").Value = 4295100
$ws.Range("D35").Value = 1516700
$ws.Range("E35").Value = 447300
$ws.Range("F35").Value = 1813000
$ws.Range("G35").Value = 2386700
$ws.Range("H35").Value = 3859600
$ws.Range("I35").Value = 4705900
$ws.Range("J35").Value = 4295100
$ws.Range("D41").Value = 1255300
$ws.Range("E41").Value = 1200900
$ws.Range("F41").Value = 4671400
$ws.Range("G41").Value = 6876000
$ws.Range("H41").Value = 4982000
$ws.Range("I41").Value = 4705200
$ws.Range("J41").Value = 4769800
$ws.Range("D42").Value = 3057700
$ws.Range("E42").Value = 2837200
$ws.Range("F42").Value = 7943000
$ws.Range("D43").Value = 10067000
$ws.Range("E43").Value = 10681000
$ws.Range("F43").Value = 22678200
$ws.Range("G43").Value = 11059400
$ws.Range("H43").Value = 10208800
$ws.Range("I43").Value = 11663400
$ws.Range("J43").Value = 12637500
$ws.Range("D44").Value = 2007200
$ws.Range("E44").Value = 1907000
$ws.Range("F44").Value = 5520200
$ws.Range("G44").Value = 3716600
$ws.Range("H44").Value = 3798200
$ws.Range("I44").Value = 2968500
$ws.Range("J44").Value = 3531600
$ws.Range("D45").Value = 1313200
$ws.Range("E45").Value = 1057500
$ws.Range("F45").Value = 4843100
$ws.Range("G45").Value = 2022200
$ws.Range("H45").Value = 1168700
$ws.Range("I45").Value = 726700
$ws.Range("J45").Value = 1084700
$ws.Range("D46").Value = 17700400
$ws.Range("E46").Value = 17683500
$ws.Range("F46").Value = 17659600
$ws.Range("G46").Value = 14916500
$ws.Range("H46").Value = 12242000
$ws.Range("I46").Value = 10816300
$ws.Range("J46").Value = 12531400
$ws.Range("D47").Value = 1332900
$ws.Range("E47").Value = 1280600
$ws.Range("F47").Value = 1364400
$ws.Range("G47").Value = 5010500
$ws.Range("H47").Value = 9194500
$ws.Range("I47").Value = 7563100
$ws.Range("J47").Value = 5608300
$ws.Range("D48").Value = 34980500
$ws.Range("E48").Value = 36265600
$ws.Range("F48").Value = 60380200
$ws.Range("G48").Value = 1843800
$ws.Range("H48").Value = 8025800
$ws.Range("I48").Value = 11282900
$ws.Range("J48").Value = 11626900
$ws.Range("D49").Value = 15257500
$ws.Range("E49").Value = 15774700
$ws.Range("F49").Value = 19995100
$ws.Range("G49").Value = 19035700
$ws.Range("H49").Value = 8736900
$ws.Range("I49").Value = 12169400
$ws.Range("J49").Value = 10804800
$ws.Range("D52").Value = 7595600
$ws.Range("E52").Value = 7353600
$ws.Range("F52").Value = 6290800
$ws.Range("G52").Value = 4880100
$ws.Range("H52").Value = 3527000
$ws.Range("I52").Value = 4861000
$ws.Range("J52").Value = 6978200
$ws.Range("D54").Value = 76866900
$ws.Range("E54").Value = 78358000
$ws.Range("F54").Value = 67054300
$ws.Range("G54").Value = 66116600
$ws.Range("H54").Value = 53043600
$ws.Range("I54").Value = 51083100
$ws.Range("J54").Value = 48596300
$ws.Range("D57").Value = 10999500
$ws.Range("E57").Value = 12271400
$ws.Range("F57").Value = 9823600
$ws.Range("G57").Value = 14950400
$ws.Range("H57").Value = 5108000
$ws.Range("I57").Value = 4747600
$ws.Range("J57").Value = 7262700
$ws.Range("D58").Value = 2676300
$ws.Range("E58").Value = 4272400
$ws.Range("F58").Value = 18555600
$ws.Range("G58").Value = 5979400
$ws.Range("H58").Value = 2673000
$ws.Range("I58").Value = 1409000
$ws.Range("J58").Value = 2756000
$ws.Range("D59").Value = 7702000
$ws.Range("E59").Value = 7764800
$ws.Range("F59").Value = 12416900
$ws.Range("G59").Value = 13947900
$ws.Range("H59").Value = 8866500
$ws.Range("I59").Value = 8744500
$ws.Range("J59").Value = 10574800
$ws.Range("D60").Value = 21377800
$ws.Range("E60").Value = 24308600
$ws.Range("F60").Value = 21992300
$ws.Range("G60").Value = 19766900
$ws.Range("H60").Value = 14169000
$ws.Range("I60").Value = 12706600
$ws.Range("J60").Value = 13363800
$ws.Range("D61").Value = 33418300
$ws.Range("E61").Value = 32335000
$ws.Range("F61").Value = 29150800
$ws.Range("G61").Value = 28236500
$ws.Range("H61").Value = 24022800
$ws.Range("I61").Value = 20897400
$ws.Range("J61").Value = 18307600
$ws.Range("D62").Value = 8590800
$ws.Range("E62").Value = 7697000
$ws.Range("F62").Value = 7763200
$ws.Range("G62").Value = 5977700
$ws.Range("H62").Value = 3975100
$ws.Range("I62").Value = 4298300
$ws.Range("J62").Value = 4695100
$ws.Range("D66").Value = 66824700
$ws.Range("E66").Value = 67552900
$ws.Range("F66").Value = 61247300
$ws.Range("G66").Value = 56580300
$ws.Range("H66").Value = 42575500
$ws.Range("I66").Value = 38381800
$ws.Range("J66").Value = 36702300
$ws.Range("D72").Value = 8848700
$ws.Range("E72").Value = 8157000
$ws.Range("F72").Value = 8913200
$ws.Range("G72").Value = 17508400
$ws.Range("H72").Value = 16551100
$ws.Range("I72").Value = 17115100
$ws.Range("J72").Value = 13901500
$ws.Range("D76").Value = 10042200
$ws.Range("E76").Value = 10805100
$ws.Range("F76").Value = 5807000
$ws.Range("G76").Value = 9536400
$ws.Range("H76").Value = 10468200
$ws.Range("I76").Value = 12701300
$ws.Range("J76").Value = 11894000
$ws.Range("D81").Value = 1516700
$ws.Range("E81").Value = 447300
$ws.Range("F81").Value = 1813000
$ws.Range("G81").Value = 2386700
$ws.Range("H81").Value = 3859600
$ws.Range("I81").Value = 4705900
$ws.Range("J81").Value = 4295100
$ws.Range("D83").Value = 8284200
$ws.Range("E83").Value = 7681800
$ws.Range("F83").Value = 6502000
$ws.Range("G83").Value = 5947500
$ws.Range("H83").Value = 5251400
$ws.Range("I83").Value = 5357400
$ws.Range("J83").Value = "NA"
$ws.Range("D89").Value = 11263200
$ws.Range("E89").Value = 12195500
$ws.Range("F89").Value = 8468000
$ws.Range("G89").Value = 12444400
$ws.Range("H89").Value = 9712500
$ws.Range("I89").Value = 10685600
$ws.Range("J89").Value = 9994400
$ws.Range("D91").Value = -6164300
$ws.Range("E91").Value = -7173900
$ws.Range("F91").Value = -6622200
$ws.Range("G91").Value = -6530400
$ws.Range("H91").Value = -6124500
$ws.Range("I91").Value = -6307600
$ws.Range("J91").Value = -6216400
$ws.Range("D94").Value = -7279300
$ws.Range("E94").Value = -7809500
$ws.Range("F94").Value = -6821900
$ws.Range("G94").Value = -7635300
$ws.Range("H94").Value = -6945600
$ws.Range("I94").Value = -10221300
$ws.Range("J94").Value = "NA"
$ws.Range("D96").Value = -832200
$ws.Range("E96").Value = -714300
$ws.Range("F96").Value = -1932200
$ws.Range("G96").Value = -882100
$ws.Range("H96").Value = -813200
$ws.Range("I96").Value = -795700
$ws.Range("J96").Value = -881500
$ws.Range("D100").Value = -3932700
$ws.Range("E100").Value = -5723100
$ws.Range("F100").Value = -2596600
$ws.Range("G100").Value = -3866500
$ws.Range("H100").Value = -2412500
$ws.Range("I100").Value = -923700
$ws.Range("J100").Value = "NA"
$ws.Range("D101").Value = 3200
$ws.Range("E101").Value = 202300
$ws.Range("F101").Value = -151800
$ws.Range("G101").Value = 4400
$ws.Range("H101").Value = -215900
$ws.Range("I101").Value = -245800
$ws.Range("J101").Value = "NA"
$ws.Range("D102").Value = 54400
$ws.Range("E102").Value = -1134800
$ws.Range("F102").Value = -1102300
$ws.Range("G102").Value = 947000
$ws.Range("H102").Value = 138400
$ws.Range("I102").Value = -705300
$ws.Range("J102").Value = -1904000
